# Fix the "Game Marvel's Spider-Man - PS4" tweets training sheet:
#  - B46 (the classification flag for that tweet) was mis-entered as 1; correct it to 0.
#  - Update the saved cursor/selection to rest on A47 (top-left, no extra scroll)
#    instead of the previous B47 selection with a scrolled viewport.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Treinamento")

# Correct the mis-labelled cell.
$ws.Range("B46").Value = 0

# Reset the view: select A47, with the window scrolled back to the top
# (clears the previous topLeftCell="A22"/selection=B47 saved state).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A47").Select()
